# Add a centered "Questionnaire 53" header (Arial, 12pt) to the document's
# only section, so printed copies can be matched back to their number.

$d = $word.ActiveDocument
$section = $d.Sections.First
$header = $section.Headers(1)
$headerRange = $header.Range

# Paragraph-level formatting for the (currently empty) header paragraph.
$headerRange.ParagraphFormat.Style = "Header"
$headerRange.ParagraphFormat.Alignment = 1

# Insert the text after the (collapsed) header range.
$headerRange.InsertAfter("Questionnaire 53")

# Apply character formatting only to the inserted text, not the trailing
# paragraph mark, so the mark's rPr is left untouched.
$textRange = $header.Range.Duplicate
$textRange.End = $textRange.End - 1
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
